# ADX-485 Regenerate data templates.
# Applies the v2.0 -> v2.1 art_inputs.xlsx template changes:
#  - bump version string
#  - replace the "year" field with "calendar_quarter" in the column
#    definitions sheet and in the data template header row
#  - drop the "art_new" field entirely (definition row + data column)
#  - add dropdown / numeric data validation to the Data Template sheet
#  - minor column width / row height touch-ups that come from the
#    template being regenerated

$wb = $excel.ActiveWorkbook

$wsDefs = $wb.Worksheets.Item("Column Definitions")
$wsData = $wb.Worksheets.Item("Data Template")

# --- Column Definitions sheet -------------------------------------------------

# 1. Bump the version/date stamp.
$wsDefs.Range("A3").Value = "version 2.1; 2020-11-10"

# 2. Turn the "year" field row into the new "calendar_quarter" field row.
#    (D12 "Required" and F12 "TRUE" are untouched -- already correct.)
$wsDefs.Range("A12").Value = "calendar_quarter"
$wsDefs.Range("B12").Value = "The calendar quarter reflected the end of reporting period. Formatted as CY20XXQY, for example CY2020Q4 for end of December 2020."
$wsDefs.Range("C12").Value = "string"
$wsDefs.Range("E12").Value = "none"

# 3. Remove the "art_new" field row entirely (row 14); everything below
#    shifts up by one row.
$wsDefs.Rows.Item(14).Delete()

# 4. Column width touch-ups from the regenerated template.
$wsDefs.Columns.Item(1).ColumnWidth = 15.88
$wsDefs.Columns.Item(2).ColumnWidth = 113.57

# 5. Row height touch-ups (regenerated template uses 15pt rows throughout,
#    rather than the old 13.8pt / 17.35pt title row).
foreach ($r in 1,3,4,7,8,9,10,11,12,13,14,19) {
    $wsDefs.Rows.Item($r).RowHeight = 15
}
$wsDefs.Range("A1").Select()

# --- Data Template sheet -------------------------------------------------------

# 1. Header row: "year" -> "calendar_quarter".
$wsData.Range("E1").Value = "calendar_quarter"

# 2. Drop the "art_new" column (G) entirely.
$wsData.Columns.Item(7).Delete()

# 3. Column width touch-ups from the regenerated template.
$wsData.Columns.Item(3).ColumnWidth = 5.14
$wsData.Columns.Item(5).ColumnWidth = 16.3

# 4. Data validation for the data-entry columns.
$wsData.Range("C2:C1001").Validation.Add(3, 1, 1, '"both,male,female"')
$wsData.Range("D2:D1001").Validation.Add(3, 1, 1, '"Y000_014,Y015_999,Y000_999"')
$wsData.Range("F2:F1001").Validation.Add(1, 1, 7, "0")

$wsData.Range("B4").Select()
